$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update column C (Förändrad) from 45180 to 45181 for existing data rows 2..19
for ($r = 2; $r -le 19; $r++) {
    $ws.Cells.Item($r, 3).Value2 = 45181
}

# Row 19 gains an explicit row height (matches author's edit)
$ws.Rows.Item(19).RowHeight = 15

# 2) Append new rows 20..23 with the new avverkningsanmälan entries
$newRows = @(
    @{ Row = 20; A = "A 42337-2023"; B = 45180; C = 45181; D = "VÄRMLANDS LÄN"; E = "HAMMARÖ"; F = "Övriga Aktiebolag"; G = 9.6 },
    @{ Row = 21; A = "A 42344-2023"; B = 45180; C = 45181; D = "VÄRMLANDS LÄN"; E = "HAMMARÖ"; F = "Övriga Aktiebolag"; G = 3 },
    @{ Row = 22; A = "A 42303-2023"; B = 45180; C = 45181; D = "VÄRMLANDS LÄN"; E = "HAMMARÖ"; F = $null;               G = 2 },
    @{ Row = 23; A = "A 42339-2023"; B = 45180; C = 45181; D = "VÄRMLANDS LÄN"; E = "HAMMARÖ"; F = "Övriga Aktiebolag"; G = 4.3 }
)

foreach ($item in $newRows) {
    $r = $item.Row

    $ws.Cells.Item($r, 1).Value2 = $item.A

    $ws.Cells.Item($r, 2).Value2 = $item.B
    $ws.Cells.Item($r, 2).NumberFormat = "YYYY-MM-DD"

    $ws.Cells.Item($r, 3).Value2 = $item.C
    $ws.Cells.Item($r, 3).NumberFormat = "YYYY-MM-DD"

    $ws.Cells.Item($r, 4).Value2 = $item.D
    $ws.Cells.Item($r, 5).Value2 = $item.E

    if ($item.F) {
        $ws.Cells.Item($r, 6).Value2 = $item.F
    }

    $ws.Cells.Item($r, 7).Value2 = $item.G

    # H..Q = 0
    for ($c = 8; $c -le 17; $c++) {
        $ws.Cells.Item($r, $c).Value2 = 0
    }

    # R column: empty, wrap-text styled cell (mirrors the rest of the table)
    $ws.Cells.Item($r, 18).WrapText = $true

    # Rows 20-22 carry an explicit row height like row 19; row 23 does not
    if ($r -le 22) {
        $ws.Rows.Item($r).RowHeight = 15
    }
}
